$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old value in B1 (its content moves into A1/A2 below)
$ws.Range("B1").ClearContents()

# Write the new data into column A
$ws.Range("A1").Value = "mngr473781"
$ws.Range("A2").Value = "vYrysEg"
